$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "7"
$ws.Range("D2").Value = -0.008800000000000002
$ws.Range("E2").Value = 0.135
$ws.Range("K2").Value = 426.08
$ws.Range("L2").Value = 0.4583378154514748
$ws.Range("M2").Value = 65.997
$ws.Range("N2").Value = 0.01293855865746549
$ws.Range("O2").Value = 0.1548934472399549
$ws.Range("P2").Value = 65.996
$ws.Range("Q2").Value = 0.0129383626097867
$ws.Range("R2").Value = 0.1548911002628614
$ws.Range("S2").Value = 0.0009999999999994458
$ws.Range("T2").Value = 0.00001515220388804712
$ws.Range("U2").Value = 3272.36
$ws.Range("V2").Value = 0.6415385821831869
$ws.Range("W2").Value = 0.1678321678321678
$ws.Range("X2").Value = 0.05592028349515488
$ws.Range("Y2").Value = 0.1119118843370129
$ws.Range("Z2").Value = 0.0969703670942061
$ws.Range("AB2").Value = 0.0535687177452608
$ws.Range("AC2").Value = -0.0535687177452608
$ws.Range("AD2").Value = 11074.3
$ws.Range("AF2").Value = 11074.3
$ws.Range("AG2").Value = 7801.939999999999
$ws.Range("AH2").Value = 0.6846510995295237
$ws.Range("AI2").Value = 0.7273187006607033
$ws.Range("AJ2").Value = 0.6046731159428153
$ws.Range("AK2").Value = 0.6526722793679688
$ws.Range("AM2").Value = -0.9450000000000001
$ws.Range("F2").ClearContents()

# Row 3
$ws.Range("D3").Value = -0.04940000000000001
$ws.Range("K3").Value = -4.07
$ws.Range("L3").Value = -0.4275210084033614
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = 0
$ws.Range("U3").Value = 4.65
$ws.Range("V3").Value = 0.07416267942583732
$ws.Range("W3").Value = -0.01911695631751996
$ws.Range("X3").Value = 0.03035135744956102
$ws.Range("Y3").Value = -0.04946831376708098
$ws.Range("Z3").Value = 0.05566925910765452
$ws.Range("AB3").Value = 0.03618051544480634
$ws.Range("AC3").Value = -0.03618051544480634
$ws.Range("AD3").Value = 22.3
$ws.Range("AF3").Value = 22.3
$ws.Range("AG3").Value = 17.65
$ws.Range("AH3").Value = 0.2623529411764706
$ws.Range("AI3").Value = 0.0859344894026975
$ws.Range("AJ3").Value = 0.2196639701306783
$ws.Range("AK3").Value = 0.06925642534824406
$ws.Range("AM3").Value = 0
$ws.Range("T3").ClearContents()
$ws.Range("AQ3").ClearContents()

# Row 4
$ws.Range("K4").Value = 19.2
$ws.Range("L4").Value = 0.08139041966935141
$ws.Range("M4").Value = 10.9
$ws.Range("N4").Value = 0.03620059780803719
$ws.Range("O4").Value = 0.5677083333333334
$ws.Range("P4").Value = 10.9
$ws.Range("Q4").Value = 0.03620059780803719
$ws.Range("R4").Value = 0.5677083333333334
$ws.Range("U4").Value = 9.42
$ws.Range("V4").Value = 0.03128528727997343
$ws.Range("W4").Value = 0.1678321678321678
$ws.Range("X4").Value = 0.03040276221716763
$ws.Range("Y4").Value = 0.1374294056150002
$ws.Range("Z4").Value = 1.166839788296978
$ws.Range("AB4").Value = 0.03761985911003193
$ws.Range("AC4").Value = -0.03761985911003193
$ws.Range("AD4").Value = 108.3
$ws.Range("AF4").Value = 108.3
$ws.Range("AG4").Value = 98.88
$ws.Range("AH4").Value = 0.264533463605276
$ws.Range("AI4").Value = 0.4520033388981636
$ws.Range("AJ4").Value = 0.2472123606180309
$ws.Range("AK4").Value = 0.4295768528977322

# Row 5
$ws.Range("D5").Value = 0.213
$ws.Range("E5").Value = 0.207
$ws.Range("K5").Value = 63.5
$ws.Range("L5").Value = 0.4191419141914192
$ws.Range("M5").Value = 12.4
$ws.Range("N5").Value = 0.05
$ws.Range("O5").Value = 0.1952755905511811
$ws.Range("P5").Value = 12.4
$ws.Range("Q5").Value = 0.05
$ws.Range("R5").Value = 0.1952755905511811
$ws.Range("U5").Value = 55.4
$ws.Range("V5").Value = 0.2233870967741935
$ws.Range("W5").Value = 0.2100562355276216
$ws.Range("X5").Value = 0.03181059329096143
$ws.Range("Y5").Value = 0.1782456422366601
$ws.Range("Z5").Value = 8.55932203389831
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.03844566871140354
$ws.Range("AC5").Value = -0.03844566871140354
$ws.Range("AD5").Value = 116.5
$ws.Range("AF5").Value = 116.5
$ws.Range("AG5").Value = 61.1
$ws.Range("AH5").Value = 0.3196159122085048
$ws.Range("AI5").Value = 0.2321642088481467
$ws.Range("AJ5").Value = 0.1976706567453898
$ws.Range("AK5").Value = 0.1368727598566308

# Row 6
$ws.Range("B6").Value = "African Export-Import Bank (MUSE:AEIB.N0004)"
$ws.Range("K6").Value = 307
$ws.Range("L6").Value = 0.7599009900990099
$ws.Range("M6").Value = 31.5
$ws.Range("N6").Value = 0.007568840405593733
$ws.Range("O6").Value = 0.1026058631921824
$ws.Range("P6").Value = 31.5
$ws.Range("Q6").Value = 0.007568840405593733
$ws.Range("R6").Value = 0.1026058631921824
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 3126.9
$ws.Range("V6").Value = 0.7513335575952713
$ws.Range("W6").Value = 0.1122322146669591
$ws.Range("X6").Value = 0.05592028349515488
$ws.Range("Y6").Value = 0.05631193117180425
$ws.Range("Z6").Value = 0.04935014169842666
$ws.Range("AB6").Value = 0.0535687177452608
$ws.Range("AC6").Value = -0.0535687177452608
$ws.Range("AD6").Value = 9800.5
$ws.Range("AF6").Value = 9800.5
$ws.Range("AG6").Value = 6673.6
$ws.Range("AH6").Value = 0.7019259004605258
$ws.Range("AI6").Value = 0.7579249383250714
$ws.Range("AJ6").Value = 0.6159071192572494
$ws.Range("AK6").Value = 0.6807156408739469
$ws.Range("F6").ClearContents()

# Row 7
$ws.Range("B7").Value = "CI Capital Holding For Financial Investments (S.A.E) (CASE:CICH)"
$ws.Range("K7").Value = 31.4
$ws.Range("L7").Value = 0.3681125439624853
$ws.Range("M7").Value = 8.920999999999999
$ws.Range("N7").Value = 0.03519132149901381
$ws.Range("O7").Value = 0.2841082802547771
$ws.Range("P7").Value = 8.92
$ws.Range("Q7").Value = 0.03518737672583826
$ws.Range("R7").Value = 0.2840764331210191
$ws.Range("S7").Value = 0.0009999999999994458
$ws.Range("T7").Value = 0.0001120950566079415
$ws.Range("U7").Value = 62.2
$ws.Range("V7").Value = 0.2453648915187377
$ws.Range("W7").Value = 0.1996185632549269
$ws.Range("X7").Value = 0.05768308505754539
$ws.Range("Y7").Value = 0.1419354781973814
$ws.Range("Z7").Value = 0.1423564753004005
$ws.Range("AB7").Value = 0.0540340203191685
$ws.Range("AC7").Value = -0.0540340203191685
$ws.Range("AD7").Value = 631.9
$ws.Range("AF7").Value = 631.9
$ws.Range("AG7").Value = 569.6999999999999
$ws.Range("AH7").Value = 0.7136887282584142
$ws.Range("AI7").Value = 0.7593126652247055
$ws.Range("AJ7").Value = 0.6920553935860058
$ws.Range("AK7").Value = 0.7398701298701298
$ws.Range("AM7").Value = 0
$ws.Range("AQ7").ClearContents()

# Row 8
$ws.Range("A8").Value = "Egypt"
$ws.Range("B8").Value = "International Company for Leasing S.A.E. (CASE:ICLE)"
$ws.Range("C8").Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Range("D8").Value = -0.207
$ws.Range("E8").Value = 0.0208
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 4.58
$ws.Range("L8").Value = 0.3982608695652174
$ws.Range("M8").Value = 1.58
$ws.Range("N8").Value = 0.03361702127659574
$ws.Range("O8").Value = 0.3449781659388647
$ws.Range("P8").Value = 1.58
$ws.Range("Q8").Value = 0.03361702127659574
$ws.Range("R8").Value = 0.3449781659388647
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 0
$ws.Range("U8").Value = 6.73
$ws.Range("V8").Value = 0.1431914893617021
$ws.Range("W8").Value = 0.1085308056872038
$ws.Range("X8").Value = 0.07521915434780686
$ws.Range("Y8").Value = 0.03331165133939692
$ws.Range("Z8").Value = 0.04642712959224869
$ws.Range("AA8").Value = 0
$ws.Range("AB8").Value = 0.06386356713763276
$ws.Range("AC8").Value = -0.06386356713763276
$ws.Range("AD8").Value = 181.6
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 181.6
$ws.Range("AG8").Value = 174.87
$ws.Range("AH8").Value = 0.7944006999125109
$ws.Range("AI8").Value = 0.8038955289951306
$ws.Range("AJ8").Value = 0.7881642403209086
$ws.Range("AK8").Value = 0.7978737965962495
$ws.Range("AL8").Value = 0
$ws.Range("AM8").Value = -0.451
$ws.Range("AQ8").Value = -0

# Row 9
$ws.Range("A9").Value = "Egypt"
$ws.Range("B9").Value = "Al Tawfeek Leasing Company (CASE:ATLC)"
$ws.Range("C9").Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Range("D9").Value = 0.0318
$ws.Range("E9").Value = 0.135
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 4.47
$ws.Range("L9").Value = 0.1401253918495298
$ws.Range("M9").Value = 0.696
$ws.Range("N9").Value = 0.02606741573033708
$ws.Range("O9").Value = 0.1557046979865772
$ws.Range("P9").Value = 0.696
$ws.Range("Q9").Value = 0.02606741573033708
$ws.Range("R9").Value = 0.1557046979865772
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 0
$ws.Range("U9").Value = 7.06
$ws.Range("V9").Value = 0.2644194756554307
$ws.Range("W9").Value = 0.2365079365079365
$ws.Range("X9").Value = 0.1279271983480231
$ws.Range("Y9").Value = 0.1085807381599134
$ws.Range("Z9").Value = 0.196356026098732
$ws.Range("AA9").Value = 0
$ws.Range("AB9").Value = 0.06550280957165021
$ws.Range("AC9").Value = -0.06550280957165021
$ws.Range("AD9").Value = 213.2
$ws.Range("AE9").Value = 0
$ws.Range("AF9").Value = 213.2
$ws.Range("AG9").Value = 206.14
$ws.Range("AH9").Value = 0.8887036265110463
$ws.Range("AI9").Value = 0.9014799154334038
$ws.Range("AJ9").Value = 0.8853289812746951
$ws.Range("AK9").Value = 0.8984483960948395
$ws.Range("AL9").Value = 0
$ws.Range("AM9").Value = -0.494
$ws.Range("AQ9").Value = -0
